$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Create the two new linked styles used to mark up code blocks:
#    "CodeSnippet" (paragraph style) and "CodeSnippet Zchn"
#    (its linked character style), based on Standard / Default
#    Paragraph Font respectively.
# ------------------------------------------------------------------
$codeStyle = $d.Styles.Add("CodeSnippet", 1)          # wdStyleTypeParagraph
$codeStyle.BaseStyle = $d.Styles.Item("Standard")
$codeStyle.QuickStyle = $true

$codeCharStyle = $d.Styles.Add("CodeSnippetZchn", 2)   # wdStyleTypeCharacter
$codeCharStyle.NameLocal = "CodeSnippet Zchn"
$codeCharStyle.BaseStyle = $d.Styles.Item("Absatz-Standardschriftart")

$codeStyle.LinkStyle = $codeCharStyle
$codeCharStyle.LinkStyle = $codeStyle

# ------------------------------------------------------------------
# 2. Apply the new "CodeSnippet" paragraph style to the paragraph
#    that holds the "[Content]" placeholder (identified via the
#    "Content" bookmark so this keeps working no matter the exact
#    paragraph index).
# ------------------------------------------------------------------
$contentBookmark = $d.Bookmarks.Item("Content")
$contentParagraph = $contentBookmark.Range.Paragraphs(1)
$contentParagraph.Style = $codeStyle

# ------------------------------------------------------------------
# 3. Re-drop the "_GoBack" bookmark six characters into the
#    placeholder text, i.e. right after "[Conte" / before "nt]".
#    Re-adding a bookmark with an existing name moves it, and in the
#    process splits the run it lands in -- exactly the
#    "[Conte" | "nt]" run split seen in the target document.
# ------------------------------------------------------------------
$splitPoint = $contentBookmark.Range.Start + 6
$goBackRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
